# fix date of email sent
# - Adds a new "Last Action Date" column (AR) with the date the email was sent,
#   recorded as text (2025-03-24) for the two "Email" leads (rows 2 and 5).
# - Fills in / corrects the "Day 1 Action 1 Complete Date" (column R) for those
#   same two rows: row 2 gets a blank (but date-formatted) completion cell,
#   row 5 gets the actual completion date 2025-03-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -------------------------------------------------
$ws.Range("AR1").Value = "Last Action Date"

# --- Row 2 (Joint CEO / Luciana Porta) ---------------------------------
# R2: give it the same date-time formatting used by the other "Complete Date"
# cells (N2) but leave the value empty, matching the source row.
$ws.Range("N2").Copy()
$ws.Range("R2").PasteSpecial(-4122)   # xlPasteFormats

# AR2: plain text "2025-03-24" (not an actual date value/format).
$ws.Range("AR2").Value = "'2025-03-24"
$ws.Range("AR2").ClearFormats()

# --- Row 5 (Vice President, Sales & Marketing / Kirk Skaufel) ---------
# R5: same date formatting as the other "Complete Date" cells, now with an
# actual completion date of 2025-03-24.
$ws.Range("N2").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value2 = 45740        # 2025-03-24

# AR5: plain text "2025-03-24".
$ws.Range("AR5").Value = "'2025-03-24"
$ws.Range("AR5").ClearFormats()

# --- Cosmetic: restore the active selection left by the editor --------
$ws.Range("I8").Select()
